$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 15055.5
$ws.Range("J3").Value = 15055.5
$ws.Range("L3").Value = 15055.5
$ws.Range("N3").Value = -15283.5
$ws.Range("H33").Value = 223.0625
$ws.Range("I33").Value = 129.5
$ws.Range("K33").Value = 129.5
$ws.Range("M33").Value = 99.5
$ws.Range("H80").Value = 737.7778
$ws.Range("J80").Value = 463.33334
$ws.Range("L80").Value = 1390.00002
$ws.Range("N80").Value = -3386.00002
$ws.Range("H83").Value = 737.7778
$ws.Range("J83").Value = 463.33334
$ws.Range("L83").Value = 4170.00006
$ws.Range("N83").Value = -14154.00006
$ws.Range("H93").Value = 39601
$ws.Range("J93").Value = 39601
$ws.Range("L93").Value = 39601
$ws.Range("N93").Value = -44593
$ws.Range("H102").Value = 15055.5
$ws.Range("J102").Value = 15055.5
$ws.Range("L102").Value = 15055.5
$ws.Range("N102").Value = -21545.5
$ws.Range("H107").Value = 2789.2
$ws.Range("I107").Value = 1069
$ws.Range("K107").Value = 1069
$ws.Range("M107").Value = 851
$ws.Range("H132").Value = 2596.5625
$ws.Range("I132").Value = 2596.5625
$ws.Range("K132").Value = 7789.6875
$ws.Range("M132").Value = -5259.6875
$ws.Range("H137").Value = 1659.0952
$ws.Range("I137").Value = 1410.5714
$ws.Range("K137").Value = 4231.7142
$ws.Range("M137").Value = -1681.7142
$ws.Range("H138").Value = 3035.6365
$ws.Range("I138").Value = 1398.6666
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 4195.9998
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = 944.0002000000004
$ws.Range("N138").Value = -25280

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 20000.25
$ws.Range("J42").Value = 20000.25
$ws.Range("L42").Value = 20000.25
$ws.Range("N42").Value = -20972.25
$ws.Range("H76").Value = 45000
$ws.Range("J76").Value = 45000
$ws.Range("L76").Value = 45000
$ws.Range("N76").Value = -45676
$ws.Range("H79").Value = 45000
$ws.Range("J79").Value = 45000
$ws.Range("L79").Value = 45000
$ws.Range("N79").Value = -47340
$ws.Range("H132").Value = 1872.2
$ws.Range("I132").Value = 1852.75
$ws.Range("J132").Value = 1950
$ws.Range("K132").Value = 5558.25
$ws.Range("L132").Value = 5850
$ws.Range("M132").Value = -3028.25
$ws.Range("N132").Value = -10910

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 11226.6
$ws.Range("J88").Value = 11226.6
$ws.Range("L88").Value = 11226.6
$ws.Range("N88").Value = -12038.6
$ws.Range("H91").Value = 11226.6
$ws.Range("J91").Value = 11226.6
$ws.Range("L91").Value = 11226.6
$ws.Range("N91").Value = -14034.6

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5540.923
$ws.Range("I31").Value = 4093.1428
$ws.Range("K31").Value = 4093.1428
$ws.Range("M31").Value = -3798.1428
$ws.Range("H34").Value = 5540.923
$ws.Range("I34").Value = 4093.1428
$ws.Range("K34").Value = 4093.1428
$ws.Range("M34").Value = -3891.1428
$ws.Range("H132").Value = 1937.7273
$ws.Range("I132").Value = 1477.8889
$ws.Range("J132").Value = 4007
$ws.Range("K132").Value = 4433.6667
$ws.Range("L132").Value = 12021
$ws.Range("M132").Value = -1903.6667
$ws.Range("N132").Value = -17081
$ws.Range("H134").Value = 3671.2856
$ws.Range("J134").Value = 3006.5
$ws.Range("L134").Value = 9019.5
$ws.Range("N134").Value = -14089.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1588.2667
$ws.Range("I131").Value = 986
$ws.Range("J131").Value = 1989.7778
$ws.Range("K131").Value = 2958
$ws.Range("L131").Value = 5969.3334
$ws.Range("M131").Value = 2082
$ws.Range("N131").Value = -16049.3334
$ws.Range("H139").Value = 2125.9092
$ws.Range("I139").Value = 2158.5
$ws.Range("K139").Value = 6475.5
$ws.Range("M139").Value = -1335.5
$ws.Range("H140").Value = 3167.7273
$ws.Range("I140").Value = 1983
$ws.Range("K140").Value = 5949
$ws.Range("M140").Value = -769

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2933.8333
$ws.Range("I132").Value = 2856.2222
$ws.Range("K132").Value = 8568.6666
$ws.Range("M132").Value = -6038.6666

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 1500
$ws.Range("I38").Value = 1500
$ws.Range("K38").Value = 1500
$ws.Range("M38").Value = -1090
$ws.Range("H62").Value = 49987.25
$ws.Range("J62").Value = 49987.25
$ws.Range("L62").Value = 49987.25
$ws.Range("N62").Value = -51235.25
$ws.Range("H65").Value = 49987.25
$ws.Range("J65").Value = 49987.25
$ws.Range("L65").Value = 149961.75
$ws.Range("N65").Value = -156201.75
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51622
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -158112
$ws.Range("H106").Value = 14999.5
$ws.Range("J106").Value = 14999.5
$ws.Range("L106").Value = 14999.5
$ws.Range("N106").Value = -17523.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H46").Value = 27500
$ws.Range("J46").Value = 27500
$ws.Range("L46").Value = 27500
$ws.Range("N46").Value = -27962
$ws.Range("H52").Value = 15000
$ws.Range("I52").Value = 15000
$ws.Range("K52").Value = 15000
$ws.Range("M52").Value = -14774
$ws.Range("H69").Value = 19653.143
$ws.Range("J69").Value = 19653.143
$ws.Range("L69").Value = 19653.143
$ws.Range("N69").Value = -21151.143
$ws.Range("H72").Value = 19653.143
$ws.Range("J72").Value = 19653.143
$ws.Range("L72").Value = 58959.429
$ws.Range("N72").Value = -66447.429
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("K81").Value = 2000
$ws.Range("M81").Value = -939
$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("K84").Value = 10000
$ws.Range("M84").Value = -4696
$ws.Range("H134").Value = 27500
$ws.Range("J134").Value = 27500
$ws.Range("L134").Value = 82500
$ws.Range("N134").Value = -87570

Write-Host "Edit complete"
